# Update NATMI LR-pair TPM-derived metrics with recomputed values (new TPM input)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 83.40125033333334
$ws.Range("H2").Value = 250.203751
$ws.Range("I2").Value = 0.9428346765536562
$ws.Range("J2").Value = 0.9428346765536562
$ws.Range("M2").Value = 12.265576
$ws.Range("N2").Value = 36.796728
$ws.Range("O2").Value = 0.3280082578429782
$ws.Range("P2").Value = 0.3280082578429782
$ws.Range("Q2").Value = 1022.964374458526
$ws.Range("R2").Value = 9206.679370126729
$ws.Range("S2").Value = 0.3092575596903126
$ws.Range("T2").Value = 0.3092575596903125

# Row 3
$ws.Range("G3").Value = 83.40125033333334
$ws.Range("H3").Value = 250.203751
$ws.Range("I3").Value = 0.9428346765536562
$ws.Range("J3").Value = 0.9428346765536562
$ws.Range("M3").Value = 7.309488999999999
$ws.Range("N3").Value = 21.928467
$ws.Range("O3").Value = 0.1954716804667316
$ws.Range("P3").Value = 0.1954716804667316
$ws.Range("Q3").Value = 609.6205218977462
$ws.Range("R3").Value = 5486.584697079717
$ws.Range("S3").Value = 0.1842974786282505
$ws.Range("T3").Value = 0.1842974786282505

# Row 4
$ws.Range("G4").Value = 83.40125033333334
$ws.Range("H4").Value = 250.203751
$ws.Range("I4").Value = 0.9428346765536562
$ws.Range("J4").Value = 0.9428346765536562
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.9467743333333334
$ws.Range("N4").Value = 2.840323
$ws.Range("O4").Value = 0.02531881092637751
$ws.Range("P4").Value = 0.02531881092637751
$ws.Range("Q4").Value = 78.96216318350812
$ws.Range("R4").Value = 710.6594686515731
$ws.Range("S4").Value = 0.02387145291049432
$ws.Range("T4").Value = 0.02387145291049432

# Row 5
$ws.Range("G5").Value = 83.40125033333334
$ws.Range("H5").Value = 250.203751
$ws.Range("I5").Value = 0.9428346765536562
$ws.Range("J5").Value = 0.9428346765536562
$ws.Range("M5").Value = 16.518665
$ws.Range("N5").Value = 49.555995
$ws.Range("O5").Value = 0.4417451352094495
$ws.Range("P5").Value = 0.4417451352094495
$ws.Range("Q5").Value = 1377.677314837472
$ws.Range("R5").Value = 12399.09583353725
$ws.Range("S5").Value = 0.4164926316743524
$ws.Range("T5").Value = 0.4164926316743524

# Row 6
$ws.Range("G6").Value = 83.40125033333334
$ws.Range("H6").Value = 250.203751
$ws.Range("I6").Value = 0.9428346765536562
$ws.Range("J6").Value = 0.9428346765536562
$ws.Range("M6").Value = 0.3536030000000001
$ws.Range("N6").Value = 1.060809
$ws.Range("O6").Value = 0.009456115554463209
$ws.Range("P6").Value = 0.009456115554463207
$ws.Range("Q6").Value = 29.49093232161767
$ws.Range("R6").Value = 265.418390894559
$ws.Range("S6").Value = 0.008915553650246317
$ws.Range("T6").Value = 0.008915553650246315

# Row 7
$ws.Range("I7").Value = 0.03022505171551549
$ws.Range("J7").Value = 0.03022505171551549
$ws.Range("M7").Value = 12.265576
$ws.Range("N7").Value = 36.796728
$ws.Range("O7").Value = 0.3280082578429782
$ws.Range("P7").Value = 0.3280082578429782
$ws.Range("Q7").Value = 32.79382047567201
$ws.Range("R7").Value = 295.144384281048
$ws.Range("S7").Value = 0.009914066556420154
$ws.Range("T7").Value = 0.009914066556420152

# Row 8
$ws.Range("I8").Value = 0.03022505171551549
$ws.Range("J8").Value = 0.03022505171551549
$ws.Range("M8").Value = 7.309488999999999
$ws.Range("N8").Value = 21.928467
$ws.Range("O8").Value = 0.1954716804667316
$ws.Range("P8").Value = 0.1954716804667316
$ws.Range("S8").Value = 0.005908141651025682
$ws.Range("T8").Value = 0.005908141651025682

# Row 9
$ws.Range("I9").Value = 0.03022505171551549
$ws.Range("J9").Value = 0.03022505171551549
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.9467743333333334
$ws.Range("N9").Value = 2.840323
$ws.Range("O9").Value = 0.02531881092637751
$ws.Range("P9").Value = 0.02531881092637751
$ws.Range("Q9").Value = 2.531340355993667
$ws.Range("R9").Value = 22.782063203943
$ws.Range("S9").Value = 0.0007652623696251189
$ws.Range("T9").Value = 0.0007652623696251188

# Row 10
$ws.Range("I10").Value = 0.03022505171551549
$ws.Range("J10").Value = 0.03022505171551549
$ws.Range("M10").Value = 16.518665
$ws.Range("N10").Value = 49.555995
$ws.Range("O10").Value = 0.4417451352094495
$ws.Range("P10").Value = 0.4417451352094495
$ws.Range("Q10").Value = 44.165079121255
$ws.Range("R10").Value = 397.485712091295
$ws.Range("S10").Value = 0.01335176955678299
$ws.Range("T10").Value = 0.01335176955678299

# Row 11
$ws.Range("I11").Value = 0.03022505171551549
$ws.Range("J11").Value = 0.03022505171551549
$ws.Range("M11").Value = 0.3536030000000001
$ws.Range("N11").Value = 1.060809
$ws.Range("O11").Value = 0.009456115554463209
$ws.Range("P11").Value = 0.009456115554463207
$ws.Range("Q11").Value = 0.9454096001410003
$ws.Range("R11").Value = 8.508686401269001
$ws.Range("S11").Value = 0.0002858115816615409
$ws.Range("T11").Value = 0.0002858115816615408

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.2062613333333333
$ws.Range("H12").Value = 0.618784
$ws.Range("I12").Value = 0.002331743669568637
$ws.Range("J12").Value = 0.002331743669568637
$ws.Range("M12").Value = 12.265576
$ws.Range("N12").Value = 36.796728
$ws.Range("O12").Value = 0.3280082578429782
$ws.Range("P12").Value = 0.3280082578429782
$ws.Range("Q12").Value = 2.529914059861333
$ws.Range("R12").Value = 22.769226538752
$ws.Range("S12").Value = 0.0007648311787916016
$ws.Range("T12").Value = 0.0007648311787916015

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.2062613333333333
$ws.Range("H13").Value = 0.618784
$ws.Range("I13").Value = 0.002331743669568637
$ws.Range("J13").Value = 0.002331743669568637
$ws.Range("M13").Value = 7.309488999999999
$ws.Range("N13").Value = 21.928467
$ws.Range("O13").Value = 0.1954716804667316
$ws.Range("P13").Value = 0.1954716804667316
$ws.Range("Q13").Value = 1.507664947125333
$ws.Range("R13").Value = 13.568984524128
$ws.Range("S13").Value = 0.0004557898535082449
$ws.Range("T13").Value = 0.0004557898535082449

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.2062613333333333
$ws.Range("H14").Value = 0.618784
$ws.Range("I14").Value = 0.002331743669568637
$ws.Range("J14").Value = 0.002331743669568637
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.9467743333333334
$ws.Range("N14").Value = 2.840323
$ws.Range("O14").Value = 0.02531881092637751
$ws.Range("P14").Value = 0.02531881092637751
$ws.Range("Q14").Value = 0.1952829363591111
$ws.Range("R14").Value = 1.757546427232
$ws.Range("S14").Value = 0.000059036977098586
$ws.Range("T14").Value = 0.000059036977098586

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.2062613333333333
$ws.Range("H15").Value = 0.618784
$ws.Range("I15").Value = 0.002331743669568637
$ws.Range("J15").Value = 0.002331743669568637
$ws.Range("M15").Value = 16.518665
$ws.Range("N15").Value = 49.555995
$ws.Range("O15").Value = 0.4417451352094495
$ws.Range("P15").Value = 0.4417451352094495
$ws.Range("Q15").Value = 3.407161867786666
$ws.Range("R15").Value = 30.66445681008
$ws.Range("S15").Value = 0.001030036422587375
$ws.Range("T15").Value = 0.001030036422587375

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.2062613333333333
$ws.Range("H16").Value = 0.618784
$ws.Range("I16").Value = 0.002331743669568637
$ws.Range("J16").Value = 0.002331743669568637
$ws.Range("M16").Value = 0.3536030000000001
$ws.Range("N16").Value = 1.060809
$ws.Range("O16").Value = 0.009456115554463209
$ws.Range("P16").Value = 0.009456115554463207
$ws.Range("Q16").Value = 0.07293462625066667
$ws.Range("R16").Value = 0.6564116362560001
$ws.Range("S16").Value = 0.00002204923758282911
$ws.Range("T16").Value = 0.0000220492375828291

# Row 17
$ws.Range("G17").Value = 1.881585
$ws.Range("H17").Value = 5.644755
$ws.Range("I17").Value = 0.02127094711161878
$ws.Range("J17").Value = 0.02127094711161878
$ws.Range("M17").Value = 12.265576
$ws.Range("N17").Value = 36.796728
$ws.Range("O17").Value = 0.3280082578429782
$ws.Range("P17").Value = 0.3280082578429782
$ws.Range("Q17").Value = 23.07872381796
$ws.Range("R17").Value = 207.70851436164
$ws.Range("S17").Value = 0.006977046304752204
$ws.Range("T17").Value = 0.006977046304752202

# Row 18
$ws.Range("G18").Value = 1.881585
$ws.Range("H18").Value = 5.644755
$ws.Range("I18").Value = 0.02127094711161878
$ws.Range("J18").Value = 0.02127094711161878
$ws.Range("M18").Value = 7.309488999999999
$ws.Range("N18").Value = 21.928467
$ws.Range("O18").Value = 0.1954716804667316
$ws.Range("P18").Value = 0.1954716804667316
$ws.Range("Q18").Value = 13.753424860065
$ws.Range("R18").Value = 123.780823740585
$ws.Range("S18").Value = 0.004157867777027094
$ws.Range("T18").Value = 0.004157867777027094

# Row 19
$ws.Range("G19").Value = 1.881585
$ws.Range("H19").Value = 5.644755
$ws.Range("I19").Value = 0.02127094711161878
$ws.Range("J19").Value = 0.02127094711161878
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 0.6666666666666666
$ws.Range("M19").Value = 0.9467743333333334
$ws.Range("N19").Value = 2.840323
$ws.Range("O19").Value = 0.02531881092637751
$ws.Range("P19").Value = 0.02531881092637751
$ws.Range("Q19").Value = 1.781436383985
$ws.Range("R19").Value = 16.032927455865
$ws.Range("S19").Value = 0.0005385550881440517
$ws.Range("T19").Value = 0.0005385550881440515

# Row 20
$ws.Range("G20").Value = 1.881585
$ws.Range("H20").Value = 5.644755
$ws.Range("I20").Value = 0.02127094711161878
$ws.Range("J20").Value = 0.02127094711161878
$ws.Range("M20").Value = 16.518665
$ws.Range("N20").Value = 49.555995
$ws.Range("O20").Value = 0.4417451352094495
$ws.Range("P20").Value = 0.4417451352094495
$ws.Range("Q20").Value = 31.081272284025
$ws.Range("R20").Value = 279.731450556225
$ws.Range("S20").Value = 0.009396337407855084
$ws.Range("T20").Value = 0.009396337407855084

# Row 21
$ws.Range("G21").Value = 1.881585
$ws.Range("H21").Value = 5.644755
$ws.Range("I21").Value = 0.02127094711161878
$ws.Range("J21").Value = 0.02127094711161878
$ws.Range("M21").Value = 0.3536030000000001
$ws.Range("N21").Value = 1.060809
$ws.Range("O21").Value = 0.009456115554463209
$ws.Range("P21").Value = 0.009456115554463207
$ws.Range("Q21").Value = 0.6653341007550001
$ws.Range("R21").Value = 5.988006906795
$ws.Range("S21").Value = 0.0002011405338403426
$ws.Range("T21").Value = 0.0002011405338403425

# Row 22
$ws.Range("E22").Value = 3
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 0.2952356666666667
$ws.Range("H22").Value = 0.885707
$ws.Range("I22").Value = 0.003337580949640955
$ws.Range("J22").Value = 0.003337580949640955
$ws.Range("M22").Value = 12.265576
$ws.Range("N22").Value = 36.796728
$ws.Range("O22").Value = 0.3280082578429782
$ws.Range("P22").Value = 0.3280082578429782
$ws.Range("Q22").Value = 3.621235507410667
$ws.Range("R22").Value = 32.591119566696
$ws.Range("S22").Value = 0.001094754112701642
$ws.Range("T22").Value = 0.001094754112701642

# Row 23
$ws.Range("E23").Value = 3
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 0.2952356666666667
$ws.Range("H23").Value = 0.885707
$ws.Range("I23").Value = 0.003337580949640955
$ws.Range("J23").Value = 0.003337580949640955
$ws.Range("M23").Value = 7.309488999999999
$ws.Range("N23").Value = 21.928467
$ws.Range("O23").Value = 0.1954716804667316
$ws.Range("P23").Value = 0.1954716804667316
$ws.Range("Q23").Value = 2.158021857907666
$ws.Range("R23").Value = 19.422196721169
$ws.Range("S23").Value = 0.0006524025569200675
$ws.Range("T23").Value = 0.0006524025569200675

# Row 24
$ws.Range("E24").Value = 3
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 0.2952356666666667
$ws.Range("H24").Value = 0.885707
$ws.Range("I24").Value = 0.003337580949640955
$ws.Range("J24").Value = 0.003337580949640955
$ws.Range("K24").Value = 2
$ws.Range("L24").Value = 0.6666666666666666
$ws.Range("M24").Value = 0.9467743333333334
$ws.Range("N24").Value = 2.840323
$ws.Range("O24").Value = 0.02531881092637751
$ws.Range("P24").Value = 0.02531881092637751
$ws.Range("Q24").Value = 0.2795215514845556
$ws.Range("R24").Value = 2.515693963361
$ws.Range("S24").Value = 0.00008450358101543885
$ws.Range("T24").Value = 0.00008450358101543885

# Row 25
$ws.Range("E25").Value = 3
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 0.2952356666666667
$ws.Range("H25").Value = 0.885707
$ws.Range("I25").Value = 0.003337580949640955
$ws.Range("J25").Value = 0.003337580949640955
$ws.Range("M25").Value = 16.518665
$ws.Range("N25").Value = 49.555995
$ws.Range("O25").Value = 0.4417451352094495
$ws.Range("P25").Value = 0.4417451352094495
$ws.Range("Q25").Value = 4.876899073718333
$ws.Range("R25").Value = 43.892091663465
$ws.Range("S25").Value = 0.001474360147871626
$ws.Range("T25").Value = 0.001474360147871627

# Row 26
$ws.Range("E26").Value = 3
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 0.2952356666666667
$ws.Range("H26").Value = 0.885707
$ws.Range("I26").Value = 0.003337580949640955
$ws.Range("J26").Value = 0.003337580949640955
$ws.Range("M26").Value = 0.3536030000000001
$ws.Range("N26").Value = 1.060809
$ws.Range("O26").Value = 0.009456115554463209
$ws.Range("P26").Value = 0.009456115554463207
$ws.Range("Q26").Value = 0.1043962174403334
$ws.Range("R26").Value = 0.9395659569630002
$ws.Range("S26").Value = 0.00003156055113217992
$ws.Range("T26").Value = 0.00003156055113217992

